$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'26.838.42"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  -1.40%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'1.873.33"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  -1.67%  "
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  -0.28%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'301.06"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  -2.14%  "
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  -0.22%  "
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.Value = "'0.5329"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  +1.30%  "
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.Value = "'0.3757"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  -1.41%  "
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'0.07177"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  -1.74%  "
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'21.64"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  -0.13%  "
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'0.8873"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  -1.96%  "
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'0.08169"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  +1.50%  "
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.Value = "'1.878.91"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  +2.52%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'93.40"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  -2.58%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'5.287"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  -1.48%  "
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  -0.22%  "
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'14.79"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  +0.34%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'0.000008550"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  -1.53%  "
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  -0.23%  "
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'26.878.23"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  -1.38%  "
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'4.984"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  -2.78%  "
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'10.68"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  -1.26%  "
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'6.393"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  -1.26%  "
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'146.41"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  -2.06%  "
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'2.276"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  -3.33%  "
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = "'1.739"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  -0.24%  "
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'18.04"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  -1.25%  "
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.Value = "'113.83"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  -2.65%  "
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.Value = "'4.731"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  -2.34%  "
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.Value = "'4.620"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  -5.74%  "
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.Value = "'0.09134"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  -1.25%  "
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.Value = "'0.8104"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  +0.85%  "
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.Value = "'0.04977"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  -2.17%  "
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.Value = "'1.176"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  -4.35%  "
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.Value = "'2.966"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  -0.24%  "
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = "'0.6047"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  +5.51%  "
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  -5.78%  "
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'2.592"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  -3.31%  "
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'0.01952"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  -1.98%  "
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'1.071"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  -1.41%  "
$c.Style = "Normal"
$c = $ws.Range("B41")
$c.Value = "'FraxShare"
$c.Style = "Normal"
$c = $ws.Range("C41")
$c.Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'6.576"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  -0.41%  "
$c.Style = "Normal"
$c = $ws.Range("B42")
$c.Value = "'Aptos"
$c.Style = "Normal"
$c = $ws.Range("C42")
$c.Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'8.876"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  -1.41%  "
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'0.5144"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'114.62"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  -1.68%  "
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'0.1493"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  -1.64%  "
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'0.9998"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  -0.28%  "
$c.Style = "Normal"
$c = $ws.Range("B47")
$c.Value = "'NEARProtocol"
$c.Style = "Normal"
$c = $ws.Range("C47")
$c.Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'1.633"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  -0.57%  "
$c.Style = "Normal"
$c = $ws.Range("B48")
$c.Value = "'EnergySwap"
$c.Style = "Normal"
$c = $ws.Range("C48")
$c.Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'9.900"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  -2.49%  "
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'37.54"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  -2.81%  "
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'0.06048"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  +1.42%  "
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'62.14"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  -3.49%  "
$c.Style = "Normal"
